$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header in H1, cloning the existing header style (bold,
# centered, bordered) from G1 so no new style entry is created.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Save values for rows 2-30 (literal data, matches the diff row-by-row)
$saveValues = @(0,0,1,1,1,0,1,0,0,0,1,0,0,0,0,0,0,0,1,0,0,0,0,0,0,0,0,1,0)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
